# MOV addition to circuit
# Adds a "Varistors (MOV)" section (rows 33-34) to the Electrical Parts List,
# suppressing voltage spikes from switching off the pumps' inductive load.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: section header "Varistors (MOV)" (bold, red font - same look
#     as the "Notes:" header style already used on the sheet) -------------
$ws.Range("A33").Value = "Varistors (MOV)"
$ws.Range("A33").Font.Bold = $true
$ws.Range("A33").Font.Color = 255

# I33 / J33 are blank placeholders carrying the plain centered style that's
# used throughout the header row of each section.
$ws.Range("I33:J33").HorizontalAlignment = -4108

# --- Row 34: the actual Varistor part-list line --------------------------
# Column order matches the header row: B=Description, C=Manufacturer/Series,
# D=Part #, E=Mouser Part #, F/G=Type, H=# Pins, I=# Connectors,
# J=Total, K=Unit Price, L=Ext Price, M=Tax, N=Shipping, O=Total $,
# P=Source, Q=Link
$ws.Range("C34").Value = "TDK"
$ws.Range("B34").Value = "Varistors 130volts 2500A"
$ws.Range("D34").Value = "B72210S0131K101"
$ws.Range("E34").Value = "871-B72210S131K101"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "NA"
$ws.Range("H34").Value = 2
$ws.Range("I34").Value = 2
$ws.Range("J34").Formula = "=H34+I34"
$ws.Range("K34").Value = 0.54
$ws.Range("L34").Formula = "=J34*K34"
$ws.Range("M34").Formula = "=ROUND(L34*6.5/100,2)"
$ws.Range("N34").Value = 0
$ws.Range("O34").Formula = "=L34+M34+N34"
$ws.Range("P34").Value = "Mouser"

# Borders around the whole new data row (thin box, like every other row)
$ws.Range("A34:Q34").Borders.LineStyle = 1
$ws.Range("A34:Q34").Borders.Weight = 2
$ws.Range("A34:Q34").Borders.ColorIndex = 1

# Centered alignment for the Type / #Pins / #Connectors / Total cells
$ws.Range("D34:J34").HorizontalAlignment = -4108

# Highlight fills: #Pins (H34) light-blue like the rest of column H,
# Total (J34) yellow to flag the new total.
$ws.Range("H34").Interior.Color = 12895428
$ws.Range("J34").Interior.Color = 65535

# Currency format for the price columns
$ws.Range("K34:O34").NumberFormat = """$""#,##0.00"

# Mouser product-page hyperlink on Q34
$url = "https://www.mouser.com/ProductDetail/EPCOS-TDK/B72210S0131K101?qs=%2fha2pyFaduhnNYlvdOReOSOGbAqbgH4hoymzM8rrfzip5Zb74YmzSA%3d%3d"
$ws.Hyperlinks.Add($ws.Range("Q34"), $url) | Out-Null

# --- Sheet view / dimension bookkeeping ----------------------------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A31").Select() | Out-Null
